$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 6).Value = 244
$ws.Cells.Item(2, 8).Value = 'bedrooms'
$ws.Cells.Item(2, 9).Value = 'target'
$ws.Cells.Item(2, 11).Value = 'j'
$ws.Cells.Item(2, 12).Value = 'stimuli/img_x0u5z.png'
$ws.Cells.Item(2, 13).Value = 92
$ws.Cells.Item(2, 14).Value = 78.16216216216216
$ws.Cells.Item(2, 15).Value = 85.08108108108108
$ws.Cells.Item(2, 16).Value = 37
$ws.Cells.Item(2, 17).Value = 10
$ws.Cells.Item(2, 18).Value = 10
$ws.Cells.Item(2, 19).Value = 10
$ws.Cells.Item(2, 20).Value = 10
$ws.Cells.Item(2, 21).Value = 10
$ws.Cells.Item(2, 22).Value = 10

# Row 3
$ws.Cells.Item(3, 3).Value = 4
$ws.Cells.Item(3, 6).Value = 245
$ws.Cells.Item(3, 8).Value = 'living_rooms'
$ws.Cells.Item(3, 9).Value = 'distractor'
$ws.Cells.Item(3, 11).Value = 'f'
$ws.Cells.Item(3, 12).Value = 'stimuli/img_emh91.png'
$ws.Cells.Item(3, 13).Value = 82.06666666666666
$ws.Cells.Item(3, 14).Value = 63.33333333333334
$ws.Cells.Item(3, 15).Value = 72.7
$ws.Cells.Item(3, 16).Value = 45
$ws.Cells.Item(3, 17).Value = 8
$ws.Cells.Item(3, 18).Value = 8
$ws.Cells.Item(3, 19).Value = 8
$ws.Cells.Item(3, 20).Value = 8
$ws.Cells.Item(3, 21).Value = 8
$ws.Cells.Item(3, 22).Value = 8

# Row 4
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 6).Value = 246
$ws.Cells.Item(4, 8).Value = 'kitchens'
$ws.Cells.Item(4, 9).Value = 'distractor'
$ws.Cells.Item(4, 11).Value = 'f'
$ws.Cells.Item(4, 12).Value = 'stimuli/img_60242.png'
$ws.Cells.Item(4, 13).Value = 78.33333333333333
$ws.Cells.Item(4, 14).Value = 57.57575757575758
$ws.Cells.Item(4, 15).Value = 67.95454545454545
$ws.Cells.Item(4, 16).Value = 33
$ws.Cells.Item(4, 17).Value = 7
$ws.Cells.Item(4, 18).Value = 7
$ws.Cells.Item(4, 19).Value = 7
$ws.Cells.Item(4, 20).Value = 7
$ws.Cells.Item(4, 21).Value = 7
$ws.Cells.Item(4, 22).Value = 7

# Row 5
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 6).Value = 247
$ws.Cells.Item(5, 8).Value = 'bedrooms'
$ws.Cells.Item(5, 9).Value = 'target'
$ws.Cells.Item(5, 11).Value = 'j'
$ws.Cells.Item(5, 12).Value = 'stimuli/img_1vq1v.png'
$ws.Cells.Item(5, 13).Value = 69.42857142857143
$ws.Cells.Item(5, 14).Value = 46.59523809523809
$ws.Cells.Item(5, 15).Value = 58.01190476190476
$ws.Cells.Item(5, 16).Value = 42
$ws.Cells.Item(5, 17).Value = 5
$ws.Cells.Item(5, 18).Value = 5
$ws.Cells.Item(5, 19).Value = 5
$ws.Cells.Item(5, 20).Value = 5
$ws.Cells.Item(5, 21).Value = 5
$ws.Cells.Item(5, 22).Value = 5

# Row 6
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 6).Value = 248
$ws.Cells.Item(6, 8).Value = 'kitchens'
$ws.Cells.Item(6, 9).Value = 'distractor'
$ws.Cells.Item(6, 11).Value = 'f'
$ws.Cells.Item(6, 12).Value = 'stimuli/img_79b5l.png'
$ws.Cells.Item(6, 13).Value = 72.74285714285715
$ws.Cells.Item(6, 14).Value = 53.31428571428572
$ws.Cells.Item(6, 15).Value = 63.02857142857143
$ws.Cells.Item(6, 16).Value = 35
$ws.Cells.Item(6, 17).Value = 6
$ws.Cells.Item(6, 18).Value = 6
$ws.Cells.Item(6, 19).Value = 6
$ws.Cells.Item(6, 20).Value = 6
$ws.Cells.Item(6, 21).Value = 6
$ws.Cells.Item(6, 22).Value = 6

# Row 7
$ws.Cells.Item(7, 3).Value = 4
$ws.Cells.Item(7, 6).Value = 249
$ws.Cells.Item(7, 8).Value = 'bedrooms'
$ws.Cells.Item(7, 9).Value = 'target'
$ws.Cells.Item(7, 11).Value = 'j'
$ws.Cells.Item(7, 12).Value = 'stimuli/img_rvssl.png'
$ws.Cells.Item(7, 13).Value = 74.25
$ws.Cells.Item(7, 14).Value = 54.33333333333334
$ws.Cells.Item(7, 15).Value = 64.29166666666667
$ws.Cells.Item(7, 16).Value = 36
$ws.Cells.Item(7, 17).Value = 6
$ws.Cells.Item(7, 18).Value = 6
$ws.Cells.Item(7, 19).Value = 6
$ws.Cells.Item(7, 20).Value = 6
$ws.Cells.Item(7, 21).Value = 6
$ws.Cells.Item(7, 22).Value = 6

# Row 8
$ws.Cells.Item(8, 3).Value = 4
$ws.Cells.Item(8, 6).Value = 250
$ws.Cells.Item(8, 8).Value = 'kitchens'
$ws.Cells.Item(8, 9).Value = 'distractor'
$ws.Cells.Item(8, 11).Value = 'f'
$ws.Cells.Item(8, 12).Value = 'stimuli/img_0mhms.png'
$ws.Cells.Item(8, 13).Value = 78
$ws.Cells.Item(8, 14).Value = 55.68571428571428
$ws.Cells.Item(8, 15).Value = 66.84285714285714
$ws.Cells.Item(8, 16).Value = 35
$ws.Cells.Item(8, 17).Value = 7
$ws.Cells.Item(8, 18).Value = 7
$ws.Cells.Item(8, 19).Value = 7
$ws.Cells.Item(8, 20).Value = 7
$ws.Cells.Item(8, 21).Value = 7
$ws.Cells.Item(8, 22).Value = 7

# Row 9
$ws.Cells.Item(9, 3).Value = 4
$ws.Cells.Item(9, 6).Value = 251
$ws.Cells.Item(9, 8).Value = 'living_rooms'
$ws.Cells.Item(9, 9).Value = 'distractor'
$ws.Cells.Item(9, 11).Value = 'f'
$ws.Cells.Item(9, 12).Value = 'stimuli/img_3m61b.png'
$ws.Cells.Item(9, 13).Value = 81.97619047619048
$ws.Cells.Item(9, 14).Value = 63.23809523809524
$ws.Cells.Item(9, 15).Value = 72.60714285714286
$ws.Cells.Item(9, 16).Value = 42
$ws.Cells.Item(9, 17).Value = 8
$ws.Cells.Item(9, 18).Value = 8
$ws.Cells.Item(9, 19).Value = 8
$ws.Cells.Item(9, 20).Value = 8
$ws.Cells.Item(9, 21).Value = 8
$ws.Cells.Item(9, 22).Value = 8

# Row 10
$ws.Cells.Item(10, 3).Value = 4
$ws.Cells.Item(10, 6).Value = 252
$ws.Cells.Item(10, 8).Value = 'kitchens'
$ws.Cells.Item(10, 9).Value = 'distractor'
$ws.Cells.Item(10, 11).Value = 'f'
$ws.Cells.Item(10, 12).Value = 'stimuli/img_e0hwx.png'
$ws.Cells.Item(10, 13).Value = 78.12121212121212
$ws.Cells.Item(10, 14).Value = 55.36363636363637
$ws.Cells.Item(10, 15).Value = 66.74242424242425
$ws.Cells.Item(10, 16).Value = 33
$ws.Cells.Item(10, 17).Value = 7
$ws.Cells.Item(10, 18).Value = 7
$ws.Cells.Item(10, 19).Value = 7
$ws.Cells.Item(10, 20).Value = 7
$ws.Cells.Item(10, 21).Value = 7
$ws.Cells.Item(10, 22).Value = 7

# Row 11
$ws.Cells.Item(11, 3).Value = 4
$ws.Cells.Item(11, 6).Value = 253
$ws.Cells.Item(11, 8).Value = 'bedrooms'
$ws.Cells.Item(11, 9).Value = 'target'
$ws.Cells.Item(11, 11).Value = 'j'
$ws.Cells.Item(11, 12).Value = 'stimuli/img_a9acb.png'
$ws.Cells.Item(11, 13).Value = 77.11428571428571
$ws.Cells.Item(11, 14).Value = 58.42857142857143
$ws.Cells.Item(11, 15).Value = 67.77142857142857
$ws.Cells.Item(11, 16).Value = 35
$ws.Cells.Item(11, 17).Value = 7
$ws.Cells.Item(11, 18).Value = 7
$ws.Cells.Item(11, 19).Value = 7
$ws.Cells.Item(11, 20).Value = 7
$ws.Cells.Item(11, 21).Value = 7
$ws.Cells.Item(11, 22).Value = 7

# Row 12
$ws.Cells.Item(12, 3).Value = 4
$ws.Cells.Item(12, 6).Value = 254
$ws.Cells.Item(12, 8).Value = 'bedrooms'
$ws.Cells.Item(12, 9).Value = 'target'
$ws.Cells.Item(12, 11).Value = 'j'
$ws.Cells.Item(12, 12).Value = 'stimuli/img_oou46.png'
$ws.Cells.Item(12, 13).Value = 75.70270270270271
$ws.Cells.Item(12, 14).Value = 54.86486486486486
$ws.Cells.Item(12, 15).Value = 65.28378378378379
$ws.Cells.Item(12, 16).Value = 37
$ws.Cells.Item(12, 17).Value = 6
$ws.Cells.Item(12, 18).Value = 6
$ws.Cells.Item(12, 19).Value = 6
$ws.Cells.Item(12, 20).Value = 6
$ws.Cells.Item(12, 21).Value = 6
$ws.Cells.Item(12, 22).Value = 6

# Row 13
$ws.Cells.Item(13, 3).Value = 4
$ws.Cells.Item(13, 6).Value = 255
$ws.Cells.Item(13, 8).Value = 'bedrooms'
$ws.Cells.Item(13, 9).Value = 'target'
$ws.Cells.Item(13, 11).Value = 'j'
$ws.Cells.Item(13, 12).Value = 'stimuli/img_2js6m.png'
$ws.Cells.Item(13, 13).Value = 40.02777777777778
$ws.Cells.Item(13, 14).Value = 20.88888888888889
$ws.Cells.Item(13, 15).Value = 30.45833333333334
$ws.Cells.Item(13, 16).Value = 36
$ws.Cells.Item(13, 17).Value = 2
$ws.Cells.Item(13, 18).Value = 2
$ws.Cells.Item(13, 19).Value = 2
$ws.Cells.Item(13, 20).Value = 2
$ws.Cells.Item(13, 21).Value = 2
$ws.Cells.Item(13, 22).Value = 2

# Row 14
$ws.Cells.Item(14, 3).Value = 4
$ws.Cells.Item(14, 6).Value = 256
$ws.Cells.Item(14, 8).Value = 'bedrooms'
$ws.Cells.Item(14, 9).Value = 'target'
$ws.Cells.Item(14, 11).Value = 'j'
$ws.Cells.Item(14, 12).Value = 'stimuli/img_juob3.png'
$ws.Cells.Item(14, 13).Value = 79.92105263157895
$ws.Cells.Item(14, 14).Value = 59.78947368421053
$ws.Cells.Item(14, 15).Value = 69.85526315789474
$ws.Cells.Item(14, 16).Value = 38
$ws.Cells.Item(14, 17).Value = 7
$ws.Cells.Item(14, 18).Value = 7
$ws.Cells.Item(14, 19).Value = 7
$ws.Cells.Item(14, 20).Value = 7
$ws.Cells.Item(14, 21).Value = 7
$ws.Cells.Item(14, 22).Value = 7

# Row 15
$ws.Cells.Item(15, 3).Value = 4
$ws.Cells.Item(15, 6).Value = 257
$ws.Cells.Item(15, 8).Value = 'living_rooms'
$ws.Cells.Item(15, 9).Value = 'distractor'
$ws.Cells.Item(15, 11).Value = 'f'
$ws.Cells.Item(15, 12).Value = 'stimuli/img_24rt2.png'
$ws.Cells.Item(15, 13).Value = 55.26829268292683
$ws.Cells.Item(15, 14).Value = 34.19512195121951
$ws.Cells.Item(15, 15).Value = 44.73170731707317
$ws.Cells.Item(15, 16).Value = 41
$ws.Cells.Item(15, 17).Value = 3
$ws.Cells.Item(15, 18).Value = 3
$ws.Cells.Item(15, 19).Value = 3
$ws.Cells.Item(15, 20).Value = 4
$ws.Cells.Item(15, 21).Value = 4
$ws.Cells.Item(15, 22).Value = 3

# Row 16
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 6).Value = 258
$ws.Cells.Item(16, 8).Value = 'bedrooms'
$ws.Cells.Item(16, 9).Value = 'target'
$ws.Cells.Item(16, 11).Value = 'j'
$ws.Cells.Item(16, 12).Value = 'stimuli/img_cogrz.png'
$ws.Cells.Item(16, 13).Value = 60.5
$ws.Cells.Item(16, 14).Value = 39.71428571428572
$ws.Cells.Item(16, 15).Value = 50.10714285714286
$ws.Cells.Item(16, 16).Value = 42
$ws.Cells.Item(16, 17).Value = 3
$ws.Cells.Item(16, 18).Value = 3
$ws.Cells.Item(16, 19).Value = 3
$ws.Cells.Item(16, 20).Value = 3
$ws.Cells.Item(16, 21).Value = 3
$ws.Cells.Item(16, 22).Value = 3

# Row 17
$ws.Cells.Item(17, 3).Value = 4
$ws.Cells.Item(17, 6).Value = 259
$ws.Cells.Item(17, 8).Value = 'living_rooms'
$ws.Cells.Item(17, 9).Value = 'distractor'
$ws.Cells.Item(17, 11).Value = 'f'
$ws.Cells.Item(17, 12).Value = 'stimuli/img_89dvt.png'
$ws.Cells.Item(17, 13).Value = 81.09756097560975
$ws.Cells.Item(17, 14).Value = 64.6829268292683
$ws.Cells.Item(17, 15).Value = 72.89024390243902
$ws.Cells.Item(17, 16).Value = 41
$ws.Cells.Item(17, 17).Value = 8
$ws.Cells.Item(17, 18).Value = 8
$ws.Cells.Item(17, 19).Value = 8
$ws.Cells.Item(17, 20).Value = 8
$ws.Cells.Item(17, 21).Value = 8
$ws.Cells.Item(17, 22).Value = 8

# Row 18
$ws.Cells.Item(18, 3).Value = 4
$ws.Cells.Item(18, 6).Value = 260
$ws.Cells.Item(18, 8).Value = 'living_rooms'
$ws.Cells.Item(18, 9).Value = 'distractor'
$ws.Cells.Item(18, 11).Value = 'f'
$ws.Cells.Item(18, 12).Value = 'stimuli/img_vh7v8.png'
$ws.Cells.Item(18, 13).Value = 78.70454545454545
$ws.Cells.Item(18, 14).Value = 59.63636363636363
$ws.Cells.Item(18, 15).Value = 69.17045454545455
$ws.Cells.Item(18, 16).Value = 44
$ws.Cells.Item(18, 17).Value = 7
$ws.Cells.Item(18, 18).Value = 7
$ws.Cells.Item(18, 19).Value = 7
$ws.Cells.Item(18, 20).Value = 7
$ws.Cells.Item(18, 21).Value = 7
$ws.Cells.Item(18, 22).Value = 7

# Row 19
$ws.Cells.Item(19, 3).Value = 4
$ws.Cells.Item(19, 6).Value = 261
$ws.Cells.Item(19, 8).Value = 'living_rooms'
$ws.Cells.Item(19, 9).Value = 'distractor'
$ws.Cells.Item(19, 11).Value = 'f'
$ws.Cells.Item(19, 12).Value = 'stimuli/img_3jnt7.png'
$ws.Cells.Item(19, 13).Value = 49.52272727272727
$ws.Cells.Item(19, 14).Value = 35.25
$ws.Cells.Item(19, 15).Value = 42.38636363636364
$ws.Cells.Item(19, 16).Value = 44
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(19, 18).Value = 3
$ws.Cells.Item(19, 19).Value = 3
$ws.Cells.Item(19, 20).Value = 3
$ws.Cells.Item(19, 21).Value = 3
$ws.Cells.Item(19, 22).Value = 4

# Row 20
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 6).Value = 262
$ws.Cells.Item(20, 8).Value = 'kitchens'
$ws.Cells.Item(20, 9).Value = 'distractor'
$ws.Cells.Item(20, 11).Value = 'f'
$ws.Cells.Item(20, 12).Value = 'stimuli/img_qihxi.png'
$ws.Cells.Item(20, 13).Value = 76.72222222222223
$ws.Cells.Item(20, 14).Value = 56.33333333333334
$ws.Cells.Item(20, 15).Value = 66.52777777777779
$ws.Cells.Item(20, 16).Value = 36
$ws.Cells.Item(20, 17).Value = 7
$ws.Cells.Item(20, 18).Value = 7
$ws.Cells.Item(20, 19).Value = 7
$ws.Cells.Item(20, 20).Value = 7
$ws.Cells.Item(20, 21).Value = 7
$ws.Cells.Item(20, 22).Value = 7

# Row 21
$ws.Cells.Item(21, 3).Value = 4
$ws.Cells.Item(21, 6).Value = 263
$ws.Cells.Item(21, 8).Value = 'living_rooms'
$ws.Cells.Item(21, 9).Value = 'distractor'
$ws.Cells.Item(21, 11).Value = 'f'
$ws.Cells.Item(21, 12).Value = 'stimuli/img_7lz7m.png'
$ws.Cells.Item(21, 13).Value = 51.5531914893617
$ws.Cells.Item(21, 14).Value = 32.87234042553192
$ws.Cells.Item(21, 15).Value = 42.21276595744681
$ws.Cells.Item(21, 16).Value = 47
$ws.Cells.Item(21, 17).Value = 3
$ws.Cells.Item(21, 18).Value = 3
$ws.Cells.Item(21, 19).Value = 3
$ws.Cells.Item(21, 20).Value = 3
$ws.Cells.Item(21, 21).Value = 3
$ws.Cells.Item(21, 22).Value = 3

# Row 22
$ws.Cells.Item(22, 3).Value = 4
$ws.Cells.Item(22, 6).Value = 264
$ws.Cells.Item(22, 8).Value = 'kitchens'
$ws.Cells.Item(22, 9).Value = 'distractor'
$ws.Cells.Item(22, 11).Value = 'f'
$ws.Cells.Item(22, 12).Value = 'stimuli/img_cxpff.png'
$ws.Cells.Item(22, 13).Value = 74.92307692307692
$ws.Cells.Item(22, 14).Value = 53.28205128205128
$ws.Cells.Item(22, 15).Value = 64.1025641025641
$ws.Cells.Item(22, 16).Value = 39
$ws.Cells.Item(22, 17).Value = 6
$ws.Cells.Item(22, 18).Value = 6
$ws.Cells.Item(22, 19).Value = 6
$ws.Cells.Item(22, 20).Value = 6
$ws.Cells.Item(22, 21).Value = 6
$ws.Cells.Item(22, 22).Value = 6

# Row 23
$ws.Cells.Item(23, 3).Value = 4
$ws.Cells.Item(23, 6).Value = 265
$ws.Cells.Item(23, 8).Value = 'kitchens'
$ws.Cells.Item(23, 9).Value = 'distractor'
$ws.Cells.Item(23, 11).Value = 'f'
$ws.Cells.Item(23, 12).Value = 'stimuli/img_xguy9.png'
$ws.Cells.Item(23, 13).Value = 78.21621621621621
$ws.Cells.Item(23, 14).Value = 57.24324324324324
$ws.Cells.Item(23, 15).Value = 67.72972972972973
$ws.Cells.Item(23, 16).Value = 37
$ws.Cells.Item(23, 17).Value = 7
$ws.Cells.Item(23, 18).Value = 7
$ws.Cells.Item(23, 19).Value = 7
$ws.Cells.Item(23, 20).Value = 7
$ws.Cells.Item(23, 21).Value = 7
$ws.Cells.Item(23, 22).Value = 7

# Row 24
$ws.Cells.Item(24, 3).Value = 4
$ws.Cells.Item(24, 6).Value = 266
$ws.Cells.Item(24, 8).Value = 'bedrooms'
$ws.Cells.Item(24, 9).Value = 'target'
$ws.Cells.Item(24, 11).Value = 'j'
$ws.Cells.Item(24, 12).Value = 'stimuli/img_t2ioc.png'
$ws.Cells.Item(24, 13).Value = 88.1891891891892
$ws.Cells.Item(24, 14).Value = 74.05405405405405
$ws.Cells.Item(24, 15).Value = 81.12162162162161
$ws.Cells.Item(24, 16).Value = 37
$ws.Cells.Item(24, 17).Value = 10
$ws.Cells.Item(24, 18).Value = 10
$ws.Cells.Item(24, 19).Value = 10
$ws.Cells.Item(24, 20).Value = 10
$ws.Cells.Item(24, 21).Value = 10
$ws.Cells.Item(24, 22).Value = 10

# Row 25
$ws.Cells.Item(25, 3).Value = 4
$ws.Cells.Item(25, 6).Value = 267
$ws.Cells.Item(25, 8).Value = 'bedrooms'
$ws.Cells.Item(25, 9).Value = 'target'
$ws.Cells.Item(25, 11).Value = 'j'
$ws.Cells.Item(25, 12).Value = 'stimuli/img_uxxo0.png'
$ws.Cells.Item(25, 13).Value = 71.74418604651163
$ws.Cells.Item(25, 14).Value = 48.44186046511628
$ws.Cells.Item(25, 15).Value = 60.09302325581395
$ws.Cells.Item(25, 16).Value = 43
$ws.Cells.Item(25, 17).Value = 5
$ws.Cells.Item(25, 18).Value = 5
$ws.Cells.Item(25, 19).Value = 5
$ws.Cells.Item(25, 20).Value = 5
$ws.Cells.Item(25, 21).Value = 5
$ws.Cells.Item(25, 22).Value = 5

# Row 26
$ws.Cells.Item(26, 3).Value = 4
$ws.Cells.Item(26, 6).Value = 268
$ws.Cells.Item(26, 8).Value = 'bedrooms'
$ws.Cells.Item(26, 9).Value = 'target'
$ws.Cells.Item(26, 11).Value = 'j'
$ws.Cells.Item(26, 12).Value = 'stimuli/img_fnu4h.png'
$ws.Cells.Item(26, 13).Value = 85.87179487179488
$ws.Cells.Item(26, 14).Value = 70.71794871794872
$ws.Cells.Item(26, 15).Value = 78.2948717948718
$ws.Cells.Item(26, 16).Value = 39
$ws.Cells.Item(26, 17).Value = 9
$ws.Cells.Item(26, 18).Value = 9
$ws.Cells.Item(26, 19).Value = 9
$ws.Cells.Item(26, 20).Value = 9
$ws.Cells.Item(26, 21).Value = 9
$ws.Cells.Item(26, 22).Value = 9

# Row 27
$ws.Cells.Item(27, 3).Value = 4
$ws.Cells.Item(27, 6).Value = 269
$ws.Cells.Item(27, 8).Value = 'bedrooms'
$ws.Cells.Item(27, 9).Value = 'target'
$ws.Cells.Item(27, 11).Value = 'j'
$ws.Cells.Item(27, 12).Value = 'stimuli/img_3h4c9.png'
$ws.Cells.Item(27, 13).Value = 85.47619047619048
$ws.Cells.Item(27, 14).Value = 67.26190476190476
$ws.Cells.Item(27, 15).Value = 76.36904761904762
$ws.Cells.Item(27, 16).Value = 42
$ws.Cells.Item(27, 17).Value = 9
$ws.Cells.Item(27, 18).Value = 9
$ws.Cells.Item(27, 19).Value = 9
$ws.Cells.Item(27, 20).Value = 9
$ws.Cells.Item(27, 21).Value = 9
$ws.Cells.Item(27, 22).Value = 9

# Row 28
$ws.Cells.Item(28, 3).Value = 4
$ws.Cells.Item(28, 6).Value = 270
$ws.Cells.Item(28, 8).Value = 'bedrooms'
$ws.Cells.Item(28, 9).Value = 'target'
$ws.Cells.Item(28, 11).Value = 'j'
$ws.Cells.Item(28, 12).Value = 'stimuli/img_cmyvx.png'
$ws.Cells.Item(28, 13).Value = 64.25
$ws.Cells.Item(28, 14).Value = 40.09375
$ws.Cells.Item(28, 15).Value = 52.171875
$ws.Cells.Item(28, 16).Value = 32
$ws.Cells.Item(28, 17).Value = 4
$ws.Cells.Item(28, 18).Value = 4
$ws.Cells.Item(28, 19).Value = 4
$ws.Cells.Item(28, 20).Value = 4
$ws.Cells.Item(28, 21).Value = 4
$ws.Cells.Item(28, 22).Value = 4

# Row 29
$ws.Cells.Item(29, 3).Value = 4
$ws.Cells.Item(29, 6).Value = 271
$ws.Cells.Item(29, 8).Value = 'bedrooms'
$ws.Cells.Item(29, 9).Value = 'target'
$ws.Cells.Item(29, 11).Value = 'j'
$ws.Cells.Item(29, 12).Value = 'stimuli/img_f4jxo.png'
$ws.Cells.Item(29, 13).Value = 82.91666666666667
$ws.Cells.Item(29, 14).Value = 65.52777777777777
$ws.Cells.Item(29, 15).Value = 74.22222222222223
$ws.Cells.Item(29, 16).Value = 36
$ws.Cells.Item(29, 17).Value = 8
$ws.Cells.Item(29, 18).Value = 8
$ws.Cells.Item(29, 19).Value = 8
$ws.Cells.Item(29, 20).Value = 8
$ws.Cells.Item(29, 21).Value = 8
$ws.Cells.Item(29, 22).Value = 8

# Row 30
$ws.Cells.Item(30, 3).Value = 4
$ws.Cells.Item(30, 6).Value = 272
$ws.Cells.Item(30, 8).Value = 'bedrooms'
$ws.Cells.Item(30, 9).Value = 'target'
$ws.Cells.Item(30, 11).Value = 'j'
$ws.Cells.Item(30, 12).Value = 'stimuli/img_72fmj.png'
$ws.Cells.Item(30, 13).Value = 53.87179487179487
$ws.Cells.Item(30, 14).Value = 36.02564102564103
$ws.Cells.Item(30, 15).Value = 44.94871794871795
$ws.Cells.Item(30, 16).Value = 39
$ws.Cells.Item(30, 17).Value = 3
$ws.Cells.Item(30, 18).Value = 3
$ws.Cells.Item(30, 19).Value = 3
$ws.Cells.Item(30, 20).Value = 3
$ws.Cells.Item(30, 21).Value = 3
$ws.Cells.Item(30, 22).Value = 3

# Row 31
$ws.Cells.Item(31, 3).Value = 4
$ws.Cells.Item(31, 6).Value = 273
$ws.Cells.Item(31, 8).Value = 'bedrooms'
$ws.Cells.Item(31, 9).Value = 'target'
$ws.Cells.Item(31, 11).Value = 'j'
$ws.Cells.Item(31, 12).Value = 'stimuli/img_wyctg.png'
$ws.Cells.Item(31, 13).Value = 33.44736842105263
$ws.Cells.Item(31, 14).Value = 11.39473684210526
$ws.Cells.Item(31, 15).Value = 22.42105263157895
$ws.Cells.Item(31, 16).Value = 38
$ws.Cells.Item(31, 17).Value = 1
$ws.Cells.Item(31, 18).Value = 1
$ws.Cells.Item(31, 19).Value = 1
$ws.Cells.Item(31, 20).Value = 1
$ws.Cells.Item(31, 21).Value = 1
$ws.Cells.Item(31, 22).Value = 1

# Row 32
$ws.Cells.Item(32, 3).Value = 4
$ws.Cells.Item(32, 6).Value = 274
$ws.Cells.Item(32, 8).Value = 'kitchens'
$ws.Cells.Item(32, 9).Value = 'distractor'
$ws.Cells.Item(32, 11).Value = 'f'
$ws.Cells.Item(32, 12).Value = 'stimuli/img_lpj57.png'
$ws.Cells.Item(32, 13).Value = 74.77777777777777
$ws.Cells.Item(32, 14).Value = 54.44444444444444
$ws.Cells.Item(32, 15).Value = 64.61111111111111
$ws.Cells.Item(32, 16).Value = 27
$ws.Cells.Item(32, 17).Value = 6
$ws.Cells.Item(32, 18).Value = 6
$ws.Cells.Item(32, 19).Value = 6
$ws.Cells.Item(32, 20).Value = 6
$ws.Cells.Item(32, 21).Value = 6
$ws.Cells.Item(32, 22).Value = 6

# Row 33
$ws.Cells.Item(33, 3).Value = 4
$ws.Cells.Item(33, 6).Value = 275
$ws.Cells.Item(33, 8).Value = 'bedrooms'
$ws.Cells.Item(33, 9).Value = 'target'
$ws.Cells.Item(33, 11).Value = 'j'
$ws.Cells.Item(33, 12).Value = 'stimuli/img_5il0t.png'
$ws.Cells.Item(33, 13).Value = 48.09523809523809
$ws.Cells.Item(33, 14).Value = 30.90476190476191
$ws.Cells.Item(33, 15).Value = 39.5
$ws.Cells.Item(33, 16).Value = 42
$ws.Cells.Item(33, 17).Value = 2
$ws.Cells.Item(33, 18).Value = 2
$ws.Cells.Item(33, 19).Value = 2
$ws.Cells.Item(33, 20).Value = 2
$ws.Cells.Item(33, 21).Value = 2
$ws.Cells.Item(33, 22).Value = 2

# Row 34
$ws.Cells.Item(34, 3).Value = 4
$ws.Cells.Item(34, 6).Value = 276
$ws.Cells.Item(34, 8).Value = 'kitchens'
$ws.Cells.Item(34, 9).Value = 'distractor'
$ws.Cells.Item(34, 11).Value = 'f'
$ws.Cells.Item(34, 12).Value = 'stimuli/img_eppte.png'
$ws.Cells.Item(34, 13).Value = 78.42424242424242
$ws.Cells.Item(34, 14).Value = 57.03030303030303
$ws.Cells.Item(34, 15).Value = 67.72727272727272
$ws.Cells.Item(34, 16).Value = 33
$ws.Cells.Item(34, 17).Value = 7
$ws.Cells.Item(34, 18).Value = 7
$ws.Cells.Item(34, 19).Value = 7
$ws.Cells.Item(34, 20).Value = 7
$ws.Cells.Item(34, 21).Value = 7
$ws.Cells.Item(34, 22).Value = 7

# Row 35
$ws.Cells.Item(35, 3).Value = 4
$ws.Cells.Item(35, 6).Value = 277
$ws.Cells.Item(35, 8).Value = 'bedrooms'
$ws.Cells.Item(35, 9).Value = 'target'
$ws.Cells.Item(35, 11).Value = 'j'
$ws.Cells.Item(35, 12).Value = 'stimuli/img_eh0no.png'
$ws.Cells.Item(35, 13).Value = 53.66666666666666
$ws.Cells.Item(35, 14).Value = 36.02564102564103
$ws.Cells.Item(35, 15).Value = 44.84615384615385
$ws.Cells.Item(35, 16).Value = 39
$ws.Cells.Item(35, 17).Value = 3
$ws.Cells.Item(35, 18).Value = 3
$ws.Cells.Item(35, 19).Value = 3
$ws.Cells.Item(35, 20).Value = 4
$ws.Cells.Item(35, 21).Value = 3
$ws.Cells.Item(35, 22).Value = 4

# Row 36
$ws.Cells.Item(36, 3).Value = 4
$ws.Cells.Item(36, 6).Value = 278
$ws.Cells.Item(36, 8).Value = 'bedrooms'
$ws.Cells.Item(36, 9).Value = 'target'
$ws.Cells.Item(36, 11).Value = 'j'
$ws.Cells.Item(36, 12).Value = 'stimuli/img_xpco9.png'
$ws.Cells.Item(36, 13).Value = 81.55555555555556
$ws.Cells.Item(36, 14).Value = 64.68888888888888
$ws.Cells.Item(36, 15).Value = 73.12222222222222
$ws.Cells.Item(36, 16).Value = 45
$ws.Cells.Item(36, 17).Value = 8
$ws.Cells.Item(36, 18).Value = 8
$ws.Cells.Item(36, 19).Value = 8
$ws.Cells.Item(36, 20).Value = 8
$ws.Cells.Item(36, 21).Value = 8
$ws.Cells.Item(36, 22).Value = 8

# Row 37
$ws.Cells.Item(37, 3).Value = 4
$ws.Cells.Item(37, 6).Value = 279
$ws.Cells.Item(37, 8).Value = 'bedrooms'
$ws.Cells.Item(37, 9).Value = 'target'
$ws.Cells.Item(37, 11).Value = 'j'
$ws.Cells.Item(37, 12).Value = 'stimuli/img_e26ut.png'
$ws.Cells.Item(37, 13).Value = 81.07692307692308
$ws.Cells.Item(37, 14).Value = 61.28205128205128
$ws.Cells.Item(37, 15).Value = 71.17948717948718
$ws.Cells.Item(37, 16).Value = 39
$ws.Cells.Item(37, 17).Value = 8
$ws.Cells.Item(37, 18).Value = 8
$ws.Cells.Item(37, 19).Value = 8
$ws.Cells.Item(37, 20).Value = 8
$ws.Cells.Item(37, 21).Value = 8
$ws.Cells.Item(37, 22).Value = 8

# Row 38
$ws.Cells.Item(38, 3).Value = 4
$ws.Cells.Item(38, 6).Value = 280
$ws.Cells.Item(38, 8).Value = 'kitchens'
$ws.Cells.Item(38, 9).Value = 'distractor'
$ws.Cells.Item(38, 11).Value = 'f'
$ws.Cells.Item(38, 12).Value = 'stimuli/img_kugyw.png'
$ws.Cells.Item(38, 13).Value = 74.25
$ws.Cells.Item(38, 14).Value = 54.10714285714285
$ws.Cells.Item(38, 15).Value = 64.17857142857143
$ws.Cells.Item(38, 16).Value = 28
$ws.Cells.Item(38, 17).Value = 6
$ws.Cells.Item(38, 18).Value = 6
$ws.Cells.Item(38, 19).Value = 6
$ws.Cells.Item(38, 20).Value = 6
$ws.Cells.Item(38, 21).Value = 6
$ws.Cells.Item(38, 22).Value = 6

# Row 39
$ws.Cells.Item(39, 3).Value = 4
$ws.Cells.Item(39, 6).Value = 281
$ws.Cells.Item(39, 8).Value = 'bedrooms'
$ws.Cells.Item(39, 9).Value = 'target'
$ws.Cells.Item(39, 11).Value = 'j'
$ws.Cells.Item(39, 12).Value = 'stimuli/img_le8uf.png'
$ws.Cells.Item(39, 13).Value = 12.88888888888889
$ws.Cells.Item(39, 14).Value = 9.222222222222221
$ws.Cells.Item(39, 15).Value = 11.05555555555556
$ws.Cells.Item(39, 16).Value = 36
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = 1
$ws.Cells.Item(39, 19).Value = 1
$ws.Cells.Item(39, 20).Value = 1
$ws.Cells.Item(39, 21).Value = 1
$ws.Cells.Item(39, 22).Value = 1

# Row 40
$ws.Cells.Item(40, 3).Value = 4
$ws.Cells.Item(40, 6).Value = 282
$ws.Cells.Item(40, 8).Value = 'bedrooms'
$ws.Cells.Item(40, 9).Value = 'target'
$ws.Cells.Item(40, 11).Value = 'j'
$ws.Cells.Item(40, 12).Value = 'stimuli/img_jp28n.png'
$ws.Cells.Item(40, 13).Value = 65.02564102564102
$ws.Cells.Item(40, 14).Value = 44.97435897435897
$ws.Cells.Item(40, 15).Value = 55
$ws.Cells.Item(40, 16).Value = 39
$ws.Cells.Item(40, 17).Value = 4
$ws.Cells.Item(40, 18).Value = 4
$ws.Cells.Item(40, 19).Value = 4
$ws.Cells.Item(40, 20).Value = 4
$ws.Cells.Item(40, 21).Value = 4
$ws.Cells.Item(40, 22).Value = 5

# Row 41
$ws.Cells.Item(41, 3).Value = 4
$ws.Cells.Item(41, 6).Value = 283
$ws.Cells.Item(41, 8).Value = 'living_rooms'
$ws.Cells.Item(41, 9).Value = 'distractor'
$ws.Cells.Item(41, 11).Value = 'f'
$ws.Cells.Item(41, 12).Value = 'stimuli/img_cehin.png'
$ws.Cells.Item(41, 13).Value = 78.86363636363636
$ws.Cells.Item(41, 14).Value = 60.02272727272727
$ws.Cells.Item(41, 15).Value = 69.44318181818181
$ws.Cells.Item(41, 16).Value = 44
$ws.Cells.Item(41, 17).Value = 7
$ws.Cells.Item(41, 18).Value = 7
$ws.Cells.Item(41, 19).Value = 7
$ws.Cells.Item(41, 20).Value = 7
$ws.Cells.Item(41, 21).Value = 7
$ws.Cells.Item(41, 22).Value = 7
